# feat: #108 fix minutes/seconds formatting in the "Общее время" (total
# time / haul) column: zero-pad single-digit minute and second components
# so e.g. "143 ч. 8 мин. 15 сек." becomes "143 ч. 08 мин. 15 сек.".
# Hours are left unpadded.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# Locate the "Общее время" column dynamically from the header row instead
# of hard-coding its index.
$timeCol = 0
for ($c = 1; $c -le $lastCol; $c++) {
    $header = $ws.Cells.Item(1, $c).Value2
    if ($header -eq "Общее время") {
        $timeCol = $c
        break
    }
}

if ($timeCol -eq 0) {
    $timeCol = 9
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $timeCol)
    $current = $cell.Value2

    if ($current -ne $null -and ($current -match '^(\d+) ч\. (\d+) мин\. (\d+) сек\.$')) {
        $hours = $matches[1]
        $minutes = $matches[2].PadLeft(2, '0')
        $seconds = $matches[3].PadLeft(2, '0')

        $newValue = $hours + " ч. " + $minutes + " мин. " + $seconds + " сек."

        if ($newValue -ne $current) {
            $cell.Value = $newValue
        }
    }
}
